$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1464435146443515
$ws.Range("C2").Value = 0.6527196652719666
$ws.Range("J2").Value = 0.01255230125523013
$ws.Range("P2").Value = 0.104602510460251
$ws.Range("S2").Value = 0.08368200836820083

# Row 3
$ws.Range("B3").Value = 0.006211180124223602
$ws.Range("C3").Value = 0.03105590062111801
$ws.Range("J3").Value = 0.03726708074534162
$ws.Range("P3").Value = 0.7142857142857143
$ws.Range("S3").Value = 0.2111801242236025

# Row 4
$ws.Range("J4").Value = 0.04
$ws.Range("P4").Value = 0.7
$ws.Range("S4").Value = 0.26

# Row 6
$ws.Range("B6").Value = 0.05454545454545454
$ws.Range("D6").Value = 0.01818181818181818
$ws.Range("E6").Value = 0.004545454545454545
$ws.Range("F6").Value = 0.08181818181818182
$ws.Range("J6").Value = 0.2545454545454545
$ws.Range("O6").Value = 0.01818181818181818
$ws.Range("Q6").Value = 0.1772727272727273
$ws.Range("R6").Value = 0.01363636363636364
$ws.Range("S6").Value = 0.3772727272727273

# Row 7
$ws.Range("B7").Value = 0.1566265060240964
$ws.Range("D7").Value = 0.03012048192771084
$ws.Range("F7").Value = 0.05421686746987952
$ws.Range("J7").Value = 0.08433734939759036
$ws.Range("O7").Value = 0.01807228915662651
$ws.Range("Q7").Value = 0.1867469879518072
$ws.Range("R7").Value = 0.05421686746987952
$ws.Range("S7").Value = 0.4156626506024096

# Row 8
$ws.Range("B8").Value = 0.08545034642032333
$ws.Range("D8").Value = 0.02771362586605081
$ws.Range("F8").Value = 0.08083140877598152
$ws.Range("J8").Value = 0.07390300230946882
$ws.Range("O8").Value = 0.009237875288683603
$ws.Range("Q8").Value = 0.1939953810623557
$ws.Range("R8").Value = 0.09468822170900693
$ws.Range("S8").Value = 0.4341801385681293

# Row 9
$ws.Range("B9").Value = 0.08695652173913043
$ws.Range("D9").Value = 0.01863354037267081
$ws.Range("F9").Value = 0.06832298136645963
$ws.Range("J9").Value = 0.1118012422360248
$ws.Range("O9").Value = 0.01863354037267081
$ws.Range("Q9").Value = 0.2795031055900621
$ws.Range("R9").Value = 0.06832298136645963
$ws.Range("S9").Value = 0.3478260869565217

# Row 10
$ws.Range("B10").Value = 0.1041666666666667
$ws.Range("D10").Value = 0.02651515151515152
$ws.Range("E10").Value = 0.001893939393939394
$ws.Range("F10").Value = 0.07102272727272728
$ws.Range("J10").Value = 0.1041666666666667
$ws.Range("O10").Value = 0.01325757575757576
$ws.Range("Q10").Value = 0.2471590909090909
$ws.Range("R10").Value = 0.06628787878787878
$ws.Range("S10").Value = 0.365530303030303

# Row 11
$ws.Range("G11").Value = 0.1376811594202899
$ws.Range("J11").Value = 0.1304347826086956
$ws.Range("K11").Value = 0.2028985507246377
$ws.Range("L11").Value = 0.5181159420289855
$ws.Range("S11").Value = 0.0108695652173913

# Row 12
$ws.Range("G12").Value = 0.762589928057554
$ws.Range("J12").Value = 0.2014388489208633
$ws.Range("K12").Value = 0.007194244604316547
$ws.Range("L12").Value = 0.007194244604316547
$ws.Range("S12").Value = 0.02158273381294964

# Row 13
$ws.Range("F13").Value = 0.02380952380952381
$ws.Range("G13").Value = 0.6904761904761905
$ws.Range("J13").Value = 0.2380952380952381
$ws.Range("S13").Value = 0.04761904761904762

# Row 15
$ws.Range("F15").Value = 0.008695652173913044
$ws.Range("H15").Value = 0.1521739130434783
$ws.Range("I15").Value = 0.1043478260869565
$ws.Range("J15").Value = 0.3652173913043478
$ws.Range("K15").Value = 0.05652173913043478
$ws.Range("M15").Value = 0.01304347826086956
$ws.Range("O15").Value = 0.05217391304347826
$ws.Range("S15").Value = 0.2478260869565218

# Row 16
$ws.Range("F16").Value = 0.02890173410404624
$ws.Range("H16").Value = 0.1791907514450867
$ws.Range("I16").Value = 0.02890173410404624
$ws.Range("J16").Value = 0.4393063583815029
$ws.Range("K16").Value = 0.1560693641618497
$ws.Range("M16").Value = 0.01734104046242774
$ws.Range("O16").Value = 0.04046242774566474
$ws.Range("S16").Value = 0.1098265895953757

# Row 17
$ws.Range("F17").Value = 0.01758241758241758
$ws.Range("H17").Value = 0.210989010989011
$ws.Range("I17").Value = 0.08131868131868132
$ws.Range("J17").Value = 0.3934065934065934
$ws.Range("K17").Value = 0.08791208791208792
$ws.Range("M17").Value = 0.02197802197802198
$ws.Range("O17").Value = 0.07472527472527472
$ws.Range("S17").Value = 0.1120879120879121

# Row 18
$ws.Range("F18").Value = 0.02255639097744361
$ws.Range("H18").Value = 0.2105263157894737
$ws.Range("I18").Value = 0.09774436090225563
$ws.Range("J18").Value = 0.3383458646616541
$ws.Range("K18").Value = 0.1052631578947368
$ws.Range("M18").Value = 0.04511278195488722
$ws.Range("O18").Value = 0.08270676691729323
$ws.Range("S18").Value = 0.09774436090225563

# Row 19
$ws.Range("F19").Value = 0.02669039145907473
$ws.Range("H19").Value = 0.2170818505338078
$ws.Range("I19").Value = 0.07384341637010676
$ws.Range("J19").Value = 0.3345195729537366
$ws.Range("K19").Value = 0.1129893238434164
$ws.Range("M19").Value = 0.0195729537366548
$ws.Range("O19").Value = 0.09608540925266904
$ws.Range("S19").Value = 0.1192170818505338
